{"js": "// The cell text \"Uno\u0161enje li\u010dnih podataka i podataka o vozilu.\" stays the\n// same, but the author clicked/edited inside it (splitting it right after\n// \"pod\"). Word tracks the last edit location with its hidden \"_GoBack\"\n// bookmark, so that bookmark moves from the end of the document (where it\n// used to sit in an otherwise-empty trailing paragraph) into the middle of\n// that sentence, splitting the run in two around an empty bookmark.\n\n// 1) Drop the \"_GoBack\" bookmark that currently sits in the trailing empty\n//    paragraph at the end of the document.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the sentence and the exact spot where it was split (\"...li\u010dnih pod\" |\n//    \"ataka i podataka o vozilu.\").\nconst body = context.document.body;\nconst prefixResults = body.search(\"Uno\u0161enje li\u010dnih pod\", { matchCase: true });\nprefixResults.load(\"items\");\nawait context.sync();\n\nif (prefixResults.items.length === 0) {\n  throw new Error('Could not find split point \"Uno\u0161enje li\u010dnih pod\" in the document.');\n}\n\n// 3) Re-insert \"_GoBack\" as a zero-length bookmark right at that split\n//    point; this naturally breaks the single run into the two runs seen in\n//    the target document.\nconst splitPoint = prefixResults.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The cell text \"Uno\u0161enje li\u010dnih podataka i podataka o vozilu.\" stays the\n# same, but the author clicked/edited inside it (splitting it right after\n# \"pod\"). Word tracks the last edit location with its hidden \"_GoBack\"\n# bookmark, so that bookmark moves from the end of the document (where it\n# used to sit in an otherwise-empty trailing paragraph) into the middle of\n# that sentence, splitting the run in two around an empty bookmark.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"_GoBack\" bookmark that currently sits in the trailing empty\n#    paragraph at the end of the document.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Find the sentence and the exact spot where it was split (\"...li\u010dnih pod\" |\n#    \"ataka i podataka o vozilu.\").\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.MatchCase = $true\n$r.Find.MatchWholeWord = $false\n$r.Find.MatchWildcards = $false\n$found = $r.Find.Execute(\"Uno\u0161enje li\u010dnih pod\")\n\nif ($found) {\n    $splitPos = $r.End\n    $caret = $d.Range($splitPos, $splitPos)\n\n    # 3) Re-insert \"_GoBack\" as a zero-length bookmark right at that split\n    #    point; this naturally breaks the single run into the two runs seen\n    #    in the target document.\n    $d.Bookmarks.Add(\"_GoBack\", $caret)\n}\n"}
